# Update two-digit multiplication answers in the table.
$d = $word.ActiveDocument

$replacements = @(
    @("79×97=7663", "45×22=990"),
    @("20×54=1080", "26×53=1378"),
    @("33×66=2178", "68×30=2040"),
    @("45×15=675",  "39×80=3120"),
    @("84×68=5712", "13×92=1196"),
    @("86×57=4902", "40×99=3960"),
    @("54×65=3510", "83×88=7304"),
    @("63×97=6111", "76×20=1520"),
    @("16×63=1008", "27×37=999"),
    @("15×49=735",  "41×70=2870"),
    @("70×88=6160", "15×22=330"),
    @("73×74=5402", "98×39=3822"),
    @("27×92=2484", "88×94=8272"),
    @("52×75=3900", "18×78=1404"),
    @("49×60=2940", "90×42=3780"),
    @("89×36=3204", "19×60=1140"),
    @("19×63=1197", "93×36=3348"),
    @("36×67=2412", "63×56=3528"),
    @("31×45=1395", "75×33=2475"),
    @("30×67=2010", "11×63=693"),
    @("28×67=1876", "58×91=5278"),
    @("85×45=3825", "62×81=5022"),
    @("66×87=5742", "42×23=966"),
    @("38×85=3230", "18×43=774"),
    @("28×34=952",  "72×15=1080")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
